$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (shortened optimizer names) for columns C:F
$ws.Range("C1").Value = "gdpa-pd"
$ws.Range("D1").Value = "gdpa-adam"
$ws.Range("E1").Value = "gdpa-noise"
$ws.Range("F1").Value = "gdpa-direct"

# Updated data values for rows 14-21, columns C:F
$data = @{
    14 = @(48, 48, 47, 47)
    15 = @(45, 46, 45, 42)
    16 = @(40, 41, 41, 39)
    17 = @(39, 39, 32, 30)
    18 = @(28, 29, 24, 23)
    19 = @(19, 20, 14, 12)
    20 = @(10, 10, 7, 7)
    21 = @(3, 4, 2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}

# Remove now-unused columns G:J (was gdpa-noise variants + gdpa-direct), shrinking the sheet's used range
$ws.Range("G1:J21").Clear()
